# Update "想去人数" (want-to-go count) figures across the four sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 42
$ws.Range("F3").Value  = 26379
$ws.Range("F4").Value  = 579
$ws.Range("F6").Value  = 595
$ws.Range("F7").Value  = 174
$ws.Range("F8").Value  = 534
$ws.Range("F10").Value = 351
$ws.Range("F12").Value = 185
$ws.Range("F16").Value = 373
$ws.Range("F17").Value = 55
$ws.Range("F18").Value = 1510
$ws.Range("F19").Value = 186
$ws.Range("F20").Value = 32
$ws.Range("F21").Value = 430

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 185

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5003
$ws.Range("F3").Value = 211

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 42
$ws.Range("F3").Value  = 5003
$ws.Range("F4").Value  = 211
$ws.Range("F5").Value  = 26379
$ws.Range("F6").Value  = 579
$ws.Range("F10").Value = 595
$ws.Range("F13").Value = 174
$ws.Range("F14").Value = 185
$ws.Range("F15").Value = 185
$ws.Range("F20").Value = 534
$ws.Range("F23").Value = 351
$ws.Range("F25").Value = 185
$ws.Range("F32").Value = 373
$ws.Range("F33").Value = 55
$ws.Range("F35").Value = 1510
$ws.Range("F36").Value = 186
$ws.Range("F38").Value = 32
$ws.Range("F39").Value = 430
